$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 70
$ws.Range("F4").Value = 9249
$ws.Range("F5").Value = 565
$ws.Range("F6").Value = 96
$ws.Range("F9").Value = 335
$ws.Range("F10").Value = 389
$ws.Range("F12").Value = 151
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 421
$ws.Range("F15").Value = 11833
$ws.Range("F21").Value = 226
$ws.Range("F25").Value = 2710
$ws.Range("F31").Value = 976
$ws.Range("F35").Value = 2611
$ws.Range("F36").Value = 3048
$ws.Range("F37").Value = 5
$ws.Range("F39").Value = 188
$ws.Range("F42").Value = 406
$ws.Range("F43").Value = 473

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 46

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 70
$ws.Range("F8").Value = 9249
$ws.Range("F9").Value = 565
$ws.Range("F13").Value = 335
$ws.Range("F14").Value = 389
$ws.Range("F15").Value = 151
$ws.Range("F16").Value = 421
$ws.Range("F17").Value = 11833
$ws.Range("F18").Value = 46
$ws.Range("F20").Value = 226
$ws.Range("F26").Value = 2710
$ws.Range("F31").Value = 5
$ws.Range("F33").Value = 976
$ws.Range("F37").Value = 2611
$ws.Range("F38").Value = 3048
$ws.Range("F39").Value = 5
$ws.Range("F41").Value = 188
$ws.Range("F43").Value = 406
$ws.Range("F44").Value = 473
